$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 74-76 (B column only; A/date values are unchanged)
$ws.Cells.Item(74, 2).Value = "['BTCUSD.SPOT']"
$ws.Cells.Item(75, 2).Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"
$ws.Cells.Item(76, 2).Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"

# Append new rows 77-87
$newRows = @(
    @{ Row = 77; Date = "2025-08-27"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 78; Date = "2025-08-28"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 79; Date = "2025-08-29"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 80; Date = "2025-08-30"; Objects = "['BTCUSD.SPOT']" },
    @{ Row = 81; Date = "2025-08-31"; Objects = "['BTCUSD.SPOT']" },
    @{ Row = 82; Date = "2025-09-01"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']" },
    @{ Row = 83; Date = "2025-09-02"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 84; Date = "2025-09-03"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 85; Date = "2025-09-04"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']" },
    @{ Row = 86; Date = "2025-09-05"; Objects = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']" },
    @{ Row = 87; Date = "2025-09-06"; Objects = "[]" }
)

$plainStyle = $ws.Cells.Item(73, 1).Style

foreach ($item in $newRows) {
    $dateCell = $ws.Cells.Item($item.Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item.Date
    $dateCell.Style = $plainStyle

    $objCell = $ws.Cells.Item($item.Row, 2)
    $objCell.NumberFormat = "@"
    $objCell.Value = $item.Objects
    $objCell.Style = $plainStyle
}
